# This script reproduces the author's edit to the Betfair Back/Lay odds workbook:
#  1) Insert one new fixture row (Colombian Primera B) above row 15, shifting the
#     four rows that followed it (15-18) down to rows 16-19.
#  2) Refresh a large number of odds values throughout the sheet (rows 3-14, plus
#     the two Brazilian Serie A rows that were pushed down to 18 and 19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15; existing rows 15-18 shift down to 16-19
$ws.Rows(15).Insert()

# --- Update changed odds values in rows 3-14 (unaffected by the row insert) ---
# Row 3
$ws.Cells.Item(3, 9).Value = 2.4
$ws.Cells.Item(3, 14).Value = 2.74
$ws.Cells.Item(3, 19).Value = 3.45
$ws.Cells.Item(3, 20).Value = 1.87
$ws.Cells.Item(3, 21).Value = 1.92
$ws.Cells.Item(3, 22).Value = 1.71
$ws.Cells.Item(3, 30).Value = 12
$ws.Cells.Item(3, 34).Value = 21
$ws.Cells.Item(3, 38).Value = 70

# Row 4
$ws.Cells.Item(4, 8).Value = 3.6
$ws.Cells.Item(4, 12).Value = 1.39
$ws.Cells.Item(4, 13).Value = 1.06
$ws.Cells.Item(4, 14).Value = 3.6
$ws.Cells.Item(4, 17).Value = 1.83
$ws.Cells.Item(4, 18).Value = 1.45
$ws.Cells.Item(4, 19).Value = 1.83
$ws.Cells.Item(4, 25).Value = 16
$ws.Cells.Item(4, 27).Value = 70
$ws.Cells.Item(4, 28).Value = 11.5
$ws.Cells.Item(4, 29).Value = 8.6
$ws.Cells.Item(4, 30).Value = 15
$ws.Cells.Item(4, 32).Value = 14
$ws.Cells.Item(4, 33).Value = 11
$ws.Cells.Item(4, 34).Value = 17
$ws.Cells.Item(4, 36).Value = 980
$ws.Cells.Item(4, 37).Value = 980
$ws.Cells.Item(4, 40).Value = 1000

# Row 5
$ws.Cells.Item(5, 6).Value = 2.08
$ws.Cells.Item(5, 7).Value = 2.34
$ws.Cells.Item(5, 10).Value = 3.15
$ws.Cells.Item(5, 11).Value = 3.65
$ws.Cells.Item(5, 13).Value = 1.09
$ws.Cells.Item(5, 14).Value = 2.96
$ws.Cells.Item(5, 17).Value = 2.18
$ws.Cells.Item(5, 18).Value = 1.25
$ws.Cells.Item(5, 19).Value = 4.1
$ws.Cells.Item(5, 20).Value = 1.89
$ws.Cells.Item(5, 21).Value = 1.89
$ws.Cells.Item(5, 23).Value = 1.75
$ws.Cells.Item(5, 24).Value = 13.5
$ws.Cells.Item(5, 25).Value = 15
$ws.Cells.Item(5, 28).Value = 9.800000000000001
$ws.Cells.Item(5, 29).Value = 9.199999999999999

# Row 7
$ws.Cells.Item(7, 8).Value = 2.36
$ws.Cells.Item(7, 9).Value = 2.66
$ws.Cells.Item(7, 11).Value = 3.6
$ws.Cells.Item(7, 13).Value = 1.09
$ws.Cells.Item(7, 17).Value = 1.97
$ws.Cells.Item(7, 22).Value = 1.6

# Row 8
$ws.Cells.Item(8, 9).Value = 4
$ws.Cells.Item(8, 10).Value = 3.25
$ws.Cells.Item(8, 12).Value = 1.34
$ws.Cells.Item(8, 18).Value = 1.32
$ws.Cells.Item(8, 19).Value = 3.1

# Row 9
$ws.Cells.Item(9, 6).Value = 5.4
$ws.Cells.Item(9, 13).Value = 1.04
$ws.Cells.Item(9, 14).Value = 3.75
$ws.Cells.Item(9, 15).Value = 1.24
$ws.Cells.Item(9, 16).Value = 2.08
$ws.Cells.Item(9, 19).Value = 2.58
$ws.Cells.Item(9, 20).Value = 1.83
$ws.Cells.Item(9, 24).Value = 24
$ws.Cells.Item(9, 25).Value = 11.5
$ws.Cells.Item(9, 27).Value = 18
$ws.Cells.Item(9, 31).Value = 19.5
$ws.Cells.Item(9, 34).Value = 27

# Row 10
$ws.Cells.Item(10, 7).Value = 2.46
$ws.Cells.Item(10, 8).Value = 2.9
$ws.Cells.Item(10, 9).Value = 3.35
$ws.Cells.Item(10, 16).Value = 2.26
$ws.Cells.Item(10, 17).Value = 1.6
$ws.Cells.Item(10, 18).Value = 1.53
$ws.Cells.Item(10, 19).Value = 2.3
$ws.Cells.Item(10, 21).Value = 2.44
$ws.Cells.Item(10, 22).Value = 1.43
$ws.Cells.Item(10, 23).Value = 1.68

# Row 11
$ws.Cells.Item(11, 6).Value = 5.3
$ws.Cells.Item(11, 9).Value = 1.7
$ws.Cells.Item(11, 28).Value = 22
$ws.Cells.Item(11, 29).Value = 9.800000000000001
$ws.Cells.Item(11, 32).Value = 44
$ws.Cells.Item(11, 34).Value = 18.5
$ws.Cells.Item(11, 37).Value = 65

# Row 12
$ws.Cells.Item(12, 20).Value = 1.8

# Row 13
$ws.Cells.Item(13, 9).Value = 4.3
$ws.Cells.Item(13, 17).Value = 1.46
$ws.Cells.Item(13, 25).Value = 27
$ws.Cells.Item(13, 26).Value = 38

# Row 14
$ws.Cells.Item(14, 6).Value = 1.76
$ws.Cells.Item(14, 15).Value = 1.26
$ws.Cells.Item(14, 20).Value = 1.77
$ws.Cells.Item(14, 23).Value = 2.28
$ws.Cells.Item(14, 33).Value = 9.6

# --- Populate the newly inserted row 15 (Colombian Primera B match) ---
# Force the League/Date/Time/Home/Away columns to Text format so that values such as
# "2025-12-02" are not auto-converted into date serial numbers by Excel.
$ws.Range("A15:E15").NumberFormat = "@"
$ws.Cells.Item(15, 1).Value = "Colombian Primera B"
$ws.Cells.Item(15, 2).Value = "2025-12-02"
$ws.Cells.Item(15, 3).Value = "17:45:00"
$ws.Cells.Item(15, 4).Value = "Cucuta Deportivo"
$ws.Cells.Item(15, 5).Value = "Real Soacha Cundinamarca FC"
$ws.Cells.Item(15, 6).Value = 1.04
$ws.Cells.Item(15, 7).Value = 1000
$ws.Cells.Item(15, 8).Value = 1.04
$ws.Cells.Item(15, 9).Value = 1000
$ws.Cells.Item(15, 10).Value = 1.02
$ws.Cells.Item(15, 11).Value = 1000
$ws.Cells.Item(15, 12).Value = 1.01
$ws.Cells.Item(15, 13).Value = 1.01
$ws.Cells.Item(15, 14).Value = 1.24
$ws.Cells.Item(15, 15).Value = 1.01
$ws.Cells.Item(15, 16).Value = 1.25
$ws.Cells.Item(15, 17).Value = 1.01
$ws.Cells.Item(15, 18).Value = 1.13
$ws.Cells.Item(15, 19).Value = 1.01
$ws.Cells.Item(15, 20).Value = 1.01
$ws.Cells.Item(15, 21).Value = 1.01
$ws.Cells.Item(15, 22).Value = 1.01
$ws.Cells.Item(15, 23).Value = 1.01
$ws.Cells.Item(15, 24).Value = 1000
$ws.Cells.Item(15, 25).Value = 1000
$ws.Cells.Item(15, 26).Value = 1000
$ws.Cells.Item(15, 27).Value = 1000
$ws.Cells.Item(15, 28).Value = 1000
$ws.Cells.Item(15, 29).Value = 1000
$ws.Cells.Item(15, 30).Value = 1000
$ws.Cells.Item(15, 31).Value = 1000
$ws.Cells.Item(15, 32).Value = 1000
$ws.Cells.Item(15, 33).Value = 1000
$ws.Cells.Item(15, 34).Value = 1000
$ws.Cells.Item(15, 35).Value = 1000
$ws.Cells.Item(15, 36).Value = 1000
$ws.Cells.Item(15, 37).Value = 1000
$ws.Cells.Item(15, 38).Value = 1000
$ws.Cells.Item(15, 39).Value = 1000
$ws.Cells.Item(15, 40).Value = 1000
$ws.Cells.Item(15, 41).Value = 1000

# --- Additional odds updates for row 18 (shifted from old row 17, Vasco Da Gama x Mirassol) ---
$ws.Cells.Item(18, 6).Value = 2.16
$ws.Cells.Item(18, 7).Value = 2.32
$ws.Cells.Item(18, 8).Value = 3.35
$ws.Cells.Item(18, 9).Value = 3.65
$ws.Cells.Item(18, 10).Value = 3.5
$ws.Cells.Item(18, 11).Value = 3.9
$ws.Cells.Item(18, 14).Value = 3.8
$ws.Cells.Item(18, 16).Value = 1.99
$ws.Cells.Item(18, 17).Value = 1.86
$ws.Cells.Item(18, 18).Value = 1.38
$ws.Cells.Item(18, 21).Value = 2.18
$ws.Cells.Item(18, 22).Value = 1.37
$ws.Cells.Item(18, 23).Value = 1.76
$ws.Cells.Item(18, 30).Value = 15
$ws.Cells.Item(18, 32).Value = 980
$ws.Cells.Item(18, 40).Value = 20

# --- Additional odds updates for row 19 (shifted from old row 18, Gremio x Fluminense) ---
$ws.Cells.Item(19, 8).Value = 2.54
$ws.Cells.Item(19, 11).Value = 3.5
$ws.Cells.Item(19, 14).Value = 2.98
$ws.Cells.Item(19, 15).Value = 1.43
$ws.Cells.Item(19, 20).Value = 1.91
$ws.Cells.Item(19, 21).Value = 1.9
